$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 430  # was 427
$ws1.Range("F8").Value = 1159  # was 1156
$ws1.Range("F9").Value = 331  # was 332
$ws1.Range("F11").Value = 871  # was 870
$ws1.Range("F12").Value = 672  # was 670
$ws1.Range("F18").Value = 2894  # was 2890
$ws1.Range("F23").Value = 312  # was 311
$ws1.Range("F24").Value = 224  # was 222
$ws1.Range("F26").Value = 5243  # was 5237
$ws1.Range("F31").Value = 295  # was 294
$ws1.Range("F32").Value = 1078  # was 1075
$ws1.Range("F34").Value = 49  # was 48
$ws1.Range("F35").Value = 284  # was 283
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1112  # was 1109
$ws2.Range("F14").Value = 603  # was 602
$ws2.Range("F15").Value = 100  # was 98
$ws2.Range("F18").Value = 2  # was 1
$ws2.Range("F25").Value = 273  # was 272
$ws2.Range("F26").Value = 3896  # was 3895
$ws2.Range("F31").Value = 48  # was 47
$ws2.Range("F33").Value = 157  # was 156
$ws2.Range("F34").Value = 32  # was 31
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2435  # was 2434
$ws3.Range("F9").Value = 1302  # was 1300
$ws3.Range("F10").Value = 352  # was 351
$ws3.Range("F11").Value = 96  # was 95
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2435  # was 2434
$ws4.Range("F7").Value = 1302  # was 1300
$ws4.Range("F8").Value = 352  # was 351
$ws4.Range("F9").Value = 96  # was 95
$ws4.Range("F11").Value = 430  # was 427
$ws4.Range("F12").Value = 807  # was 806
$ws4.Range("F14").Value = 1159  # was 1156
$ws4.Range("F15").Value = 331  # was 332
$ws4.Range("F16").Value = 871  # was 870
$ws4.Range("F17").Value = 672  # was 670
$ws4.Range("F18").Value = 1112  # was 1109
$ws4.Range("F19").Value = 1112  # was 1109
$ws4.Range("F24").Value = 2894  # was 2890
$ws4.Range("F28").Value = 312  # was 311
$ws4.Range("F29").Value = 224  # was 222
$ws4.Range("F30").Value = 5244  # was 5237
$ws4.Range("F33").Value = 603  # was 602
$ws4.Range("F34").Value = 603  # was 602
$ws4.Range("F36").Value = 100  # was 98
$ws4.Range("F38").Value = 295  # was 294
$ws4.Range("F45").Value = 273  # was 272
$ws4.Range("F46").Value = 1078  # was 1075
$ws4.Range("F48").Value = 48  # was 47
$ws4.Range("F49").Value = 157  # was 156
$ws4.Range("F50").Value = 49  # was 48
$ws4.Range("F51").Value = 284  # was 283
